$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new ranking row (row 6) for a second player
$ws.Range("A6").Value = "Tomek1"
$ws.Range("B6").Value = "03:42"
$ws.Range("C6").Value = 84
$ws.Range("D6").Value = "Galactic Tower"

# Update the view: zoom + selected cell
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("E14").Select()
